{"js": "// Rename four of the document's auto-generated bookmarks to new\n// (regenerated) random names, leaving their position/content untouched.\n// Office.js has no \"rename\" verb for a bookmark, so each rename is done by\n// capturing the bookmark's (empty, zero-length) range, deleting the old\n// bookmark, and inserting a new bookmark with the new name at that same\n// range.\nconst renames = [\n  [\"_o0mfdio1out9\", \"_7hmavikhyc68\"],\n  [\"_k5n8ulqad0l3\", \"_gjx167ng6cvb\"],\n  [\"_mvnfuwyjrtfl\", \"_pmxekio815o8\"],\n  [\"_sujk5djs6qjr\", \"_6rn4t82o259w\"],\n];\n\nconst doc = context.document;\n\nfor (const [oldName, newName] of renames) {\n  const range = doc.getBookmarkRange(oldName);\n  doc.deleteBookmark(oldName);\n  range.insertBookmark(newName);\n}\n\nawait context.sync();\n", "ps1": "# Rename four of the document's auto-generated bookmarks to new\n# (regenerated) random names, leaving their position/content untouched.\n# The Word object model's Bookmark object has no in-place rename, so each\n# rename is done by grabbing the bookmark's Range, deleting the bookmark,\n# and adding a new bookmark with the new name over that same Range.\n\n$d = $word.ActiveDocument\n\n$renames = @(\n    @(\"_o0mfdio1out9\", \"_7hmavikhyc68\"),\n    @(\"_k5n8ulqad0l3\", \"_gjx167ng6cvb\"),\n    @(\"_mvnfuwyjrtfl\", \"_pmxekio815o8\"),\n    @(\"_sujk5djs6qjr\", \"_6rn4t82o259w\")\n)\n\nforeach ($pair in $renames) {\n    $oldName = $pair[0]\n    $newName = $pair[1]\n    $rng = $d.Bookmarks($oldName).Range\n    $d.Bookmarks($oldName).Delete()\n    $d.Bookmarks.Add($newName, $rng)\n}\n"}
